# "update existing employee job"
# The employees in rows 2-11 are relocated to a new Chicago office: their
# EmpId, Address and City (and the ZipCode in column E) are updated, while
# State (col D), DisplayName (col F) and EffectiveDate (col G) stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  EmpId = "234563"; Address = "123 Main St";         Zip = 60176 },
    @{ Row = 3;  EmpId = "002031"; Address = "234 Main St";         Zip = 60608 },
    @{ Row = 4;  EmpId = "002542"; Address = "234 Durham Ln";       Zip = 60612 },
    @{ Row = 5;  EmpId = "976537"; Address = "325 MLK BLVD";        Zip = 60609 },
    @{ Row = 6;  EmpId = "000727"; Address = "455 Dearborn Ave";    Zip = 60618 },
    @{ Row = 7;  EmpId = "000065"; Address = "3114 Rudder Ln";      Zip = 60610 },
    @{ Row = 8;  EmpId = "002595"; Address = "987 Express Parkway"; Zip = 60614 },
    @{ Row = 9;  EmpId = "002756"; Address = "5443 Glenbridge Rd";  Zip = 60617 },
    @{ Row = 10; EmpId = "000002"; Address = "5678 Hemingway St";   Zip = 60618 },
    @{ Row = 11; EmpId = "002452"; Address = "67677 Lord St";       Zip = 60611 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.EmpId      # A: EmpId
    $ws.Cells.Item($r, 2).Value = $u.Address    # B: Adderss
    $ws.Cells.Item($r, 3).Value = "Chicago"     # C: City
    $ws.Cells.Item($r, 5).Value = $u.Zip        # E: ZipCode
}

$ws.Range("A11").Select()
